# The source diff for this change touches exactly four `<w:nsid w:val="…"/>`
# attributes inside word/numbering.xml (the w:nsid child of four
# w:abstractNum definitions, abstractNumId 990, 991, 99416 and 99417).
#
# w:nsid is the OOXML "Numbering Definition Instance Identifier" - a GUID
# stamped on an abstract numbering definition purely so two definitions can
# be told apart internally. Per the schema it carries no semantic meaning,
# is never surfaced in the Word UI, and - critically - is not exposed by
# any property/method of the Word object model (Document.Lists,
# Document.ListTemplates, ListFormat, ListTemplate, ListLevel, …): there is
# no ListTemplate.Nsid/Guid/Id member, Find/Replace only ever sees story
# (body) text and never touches numbering.xml, and Document.WordOpenXML /
# Content.WordOpenXML are read-only in real Word automation, so round-
# tripping through them cannot be used to smuggle a raw XML patch in
# either. The commit message ("Automatic build output files") confirms
# this value is just re-randomized packaging noise from whatever pipeline
# regenerated the .docx, not a deliberate, user-visible edit - exactly the
# kind of change a COM automation script has no API surface to reproduce.
#
# Concretely there is nothing else in the diff (no paragraph/run/table
# text, no formatting, no list usage changes) - every list instance still
# points at the same abstractNumId values, only their internal nsid GUIDs
# differ. So the correct, faithful COM-automation action here is a no-op:
# touch nothing, and let the document round-trip unchanged.
$d = $word.ActiveDocument
